$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextCell "D2" "95.404.45"
Set-TextCell "E2" "  -2.27%  "

# Row 3 (Ethereum)
Set-TextCell "D3" "3.616.29"
Set-TextCell "E3" "  -2.91%  "

# Row 4 (XRP)
Set-TextCell "E4" "  +26.09%  "

# Row 5 (TetherUSD)
Set-TextCell "E5" "  +0.04%  "

# Row 6 (Solana)
Set-TextCell "D6" "223.70"
Set-TextCell "E6" "  -6.18%  "

# Row 7 (BNB)
Set-TextCell "D7" "641.11"
Set-TextCell "E7" "  -2.58%  "

# Row 8 (Dogecoin)
Set-TextCell "E8" "  -5.02%  "

# Row 9 (Cardano)
Set-TextCell "D9" "1.19"
Set-TextCell "E9" "  +5.48%  "

# Row 10 (USDC)
Set-TextCell "D10" "1.00"
Set-TextCell "E10" "  -0.01%  "

# Row 11 (LidoStakedEther)
Set-TextCell "D11" "3.613.09"
Set-TextCell "E11" "  -2.93%  "

# Row 12 (Avalanche)
Set-TextCell "D12" "50.38"
Set-TextCell "E12" "  +12.24%  "

# Row 13 (TRON)
Set-TextCell "E13" "  +4.83%  "

# Row 14 (ShibaInu)
Set-TextCell "D14" "0.0000294"
Set-TextCell "E14" "  -6.59%  "

# Row 15 (Toncoin)
Set-TextCell "D15" "6.51"
Set-TextCell "E15" "  -5.01%  "

# Row 16 (WrappedliquidstakedEther2.0)
Set-TextCell "D16" "4.288.70"
Set-TextCell "E16" "  -3.09%  "

# Row 17 (WrappedBTC)
Set-TextCell "D17" "95.300.44"
Set-TextCell "E17" "  -2.08%  "

# Row 18 (Chainlink)
Set-TextCell "D18" "24.38"
Set-TextCell "E18" "  +29.25%  "

# Row 19 (Polkadot)
Set-TextCell "D19" "8.99"
Set-TextCell "E19" "  -2.85%  "

# Row 20 (Uniswap)
Set-TextCell "E20" "  +4.64%  "

# Row 21 (WrappedEther)
Set-TextCell "D21" "3.612.31"
Set-TextCell "E21" "  -2.99%  "

# Row 22 (Hedera)
Set-TextCell "D22" "0.292"
Set-TextCell "E22" "  +37.45%  "

# Row 23 (Stellar)
Set-TextCell "D23" "0.535"
Set-TextCell "E23" "  -0.98%  "

# Row 24 (Litecoin)
Set-TextCell "D24" "136.63"
Set-TextCell "E24" "  +15.64%  "

# Row 25 (BitcoinCash)
Set-TextCell "D25" "532.43"
Set-TextCell "E25" "  +0.49%  "

# Row 26 (SuiNetwork)
Set-TextCell "D26" "3.28"
Set-TextCell "E26" "  -5.67%  "

# Row 27 (NEARProtocol)
Set-TextCell "D27" "7.02"
Set-TextCell "E27" "  +1.59%  "

# Row 28 (PEPE)
Set-TextCell "D28" "0.0000203"
Set-TextCell "E28" "  -9.74%  "

# Row 29 (Aptos)
Set-TextCell "D29" "13.21"
Set-TextCell "E29" "  -1.74%  "

# Row 30 (WrappedeETH)
Set-TextCell "D30" "3.784.17"
Set-TextCell "E30" "  -3.63%  "

# Row 31 (InternetComputer(DFINITY))
Set-TextCell "D31" "13.39"
Set-TextCell "E31" "  +4.49%  "

# Row 32 (PancakeSwap)
Set-TextCell "D32" "3.15"
Set-TextCell "E32" "  +3.40%  "

# Row 33 (Dai)
Set-TextCell "E33" "  -0.08%  "

# Row 34: was Fetch.AI, now PolygonEcosystemToken
$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell "D34" "0.642"
Set-TextCell "E34" "  +7.12%  "

# Row 35: was PolygonEcosystemToken, now Fetch.AI
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D35" "1.87"
Set-TextCell "E35" "  +1.60%  "

# Row 36 (EthereumClassic)
Set-TextCell "D36" "33.67"
Set-TextCell "E36" "  +1.54%  "

# Row 37 (Cronos)
Set-TextCell "E37" "  -4.04%  "

# Row 38 (Binance-PegBSC-USD)
Set-TextCell "E38" "  +0.04%  "

# Row 39 (VeChain)
Set-TextCell "D39" "0.0558"
Set-TextCell "E39" "  +21.52%  "

# Row 40: was USDe, now Filecoin
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D40" "7.35"
Set-TextCell "E40" "  +7.53%  "

# Row 41: was RenderToken, now USDe
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D41" "1.00"
Set-TextCell "E41" "  -0.04%  "

# Row 42: was Filecoin, now RenderToken
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell "D42" "8.57"
Set-TextCell "E42" "  -2.32%  "

# Row 43 (Bittensor)
Set-TextCell "D43" "591.03"
Set-TextCell "E43" "  -7.95%  "

# Row 44 (Algorand)
Set-TextCell "D44" "0.508"
Set-TextCell "E44" "  +2.55%  "

# Row 45 (ARBITRUM)
Set-TextCell "D45" "1.01"
Set-TextCell "E45" "  +4.69%  "

# Row 46 (EnergySwap)
Set-TextCell "D46" "40.90"
Set-TextCell "E46" "  -0.75%  "

# Row 47 (ImmutableX)
Set-TextCell "D47" "2.01"
Set-TextCell "E47" "  -0.04%  "

# Row 48 (Kaspa)
Set-TextCell "D48" "0.157"
Set-TextCell "E48" "  -7.02%  "

# Row 49 (Cosmos)
Set-TextCell "D49" "9.33"
Set-TextCell "E49" "  +5.95%  "

# Row 50 (Aave)
Set-TextCell "D50" "234.60"
Set-TextCell "E50" "  +12.12%  "

# Row 51 (Stacks)
Set-TextCell "D51" "2.34"
Set-TextCell "E51" "  -2.60%  "
